# The sheet previously used row 1 as the text header row
# (Lg.,mm / Threading / HeadDia., mm / ... / material_surface) and row 2
# started the first data group ("M2 x 0.4 mm").
#
# The new layout inserts a brand-new row 1 containing a simple numeric
# index sequence (0,1,2,...,11) across columns A:L, carrying the bold /
# bordered / centered header style that used to sit on the text header
# row. Everything that used to be in row 1 (and below) shifts down by
# one row, and the former header row (now row 2) loses that bold style,
# becoming a plain row like the rest of the data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at the very top; this shifts every existing
# row (old header row + all data rows) down by one.
$ws.Rows.Item(1).Insert()

# Populate the new row 1 with the numeric sequence 0..11 across A1:L1.
for ($col = 1; $col -le 12; $col++) {
    $ws.Cells.Item(1, $col).Value = $col - 1
}

# Give the new row 1 the same formatting the old header row (now row 2)
# has (bold font, thin border, centered alignment) by copying its
# formatting over.
$ws.Range("A2:L2").Copy()
$ws.Range("A1:L1").PasteSpecial(-4122)

# The old header row (now row 2) should no longer carry that special
# formatting - reset it back to the default Normal style.
$ws.Range("A2:L2").Style = "Normal"
